$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transformer_types")
$ws.Activate()

$ws.Cells.Item(1, 3).Value = "HV"
$ws.Cells.Item(1, 4).Value = "LV"
$ws.Cells.Item(1, 5).Value = "rating"
$ws.Cells.Item(1, 6).Value = "Pcu"
$ws.Cells.Item(1, 7).Value = "Pfe"
$ws.Cells.Item(1, 8).Value = "I0"
$ws.Cells.Item(1, 9).Value = "Vsc"
